$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: the SmartScore cells were entered as plain text ("0.563" etc.)
# and get normalized to real numbers (Excel's "convert to number" behaviour).
$ws.Range("I3").Value = 0.563
$ws.Range("L3").Value = 0.454
$ws.Range("O3").Value = 0.418
$ws.Range("R3").Value = 0.663
$ws.Range("U3").Value = 0.608
$ws.Range("X3").Value = 0.576
$ws.Range("AA3").Value = 0.730
$ws.Range("AD3").Value = 0.580
$ws.Range("AG3").Value = 0.556

# --- Row 4: new submission appended by the Streamlit app
# (Julieta Hernandez, 2025-11-13 20:43:29)
$ws.Range("A4").Value = "Julieta Hernandez_20251113_204328"

# B4 stays blank text (matches the empty Grupo_Experimental column on prior rows).
$ws.Range("B4").Formula = '=""'

$ws.Range("C4").Value = "Julieta Hernandez"
$ws.Range("D4").Value = 22
$ws.Range("E4").Value = "Female"
$ws.Range("F4").Value = "2025-11-13 20:43:29"
$ws.Range("G4").Value = "{
  ""portion"": 0.6,
  ""diet"": 0.14285714285714285,
  ""salt"": 0.4,
  ""fat"": 1.0,
  ""natural"": 0.4,
  ""convenience"": 0.6,
  ""price"": 0.8
}"

$ws.Range("H4").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("I4").Value = "'0.612"
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"

$ws.Range("K4").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("L4").Value = "'0.573"
$ws.Range("L4").Style = "Normal"
$ws.Range("M4").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

$ws.Range("N4").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("O4").Value = "'0.496"
$ws.Range("O4").Style = "Normal"
$ws.Range("P4").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

$ws.Range("Q4").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("R4").Value = "'0.643"
$ws.Range("R4").Style = "Normal"
$ws.Range("S4").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

$ws.Range("T4").Value = "Velveeta Original Shells & Cheese (microwave cups)"
$ws.Range("U4").Value = "'0.626"
$ws.Range("U4").Style = "Normal"
$ws.Range("V4").Value = "Muy cremoso, porción individual, rápido, salado, ideal para niños"

$ws.Range("W4").Value = "Annie’s Shells & White Cheddar"
$ws.Range("X4").Value = "'0.568"
$ws.Range("X4").Style = "Normal"
$ws.Range("Y4").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"

$ws.Range("Z4").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AA4").Value = "'0.695"
$ws.Range("AA4").Style = "Normal"
$ws.Range("AB4").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

$ws.Range("AC4").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AD4").Value = "'0.686"
$ws.Range("AD4").Style = "Normal"
$ws.Range("AE4").Value = "Portátil, saludable, fácil, buena textura, sabor suave"

$ws.Range("AF4").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AG4").Value = "'0.645"
$ws.Range("AG4").Style = "Normal"
$ws.Range("AH4").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

# The multi-line JSON in G4 otherwise triggers an autosized custom row height;
# AutoFit snaps row 4 back to the sheet's normal (default) height, same as rows 2-3.
$ws.Rows(4).AutoFit()
